$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.803.44"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "1.866.63"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.52"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07787"
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3086"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.87"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07854"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.193"
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").Value = "1.864.51"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.98"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6975"
$ws.Range("E15").Value = "  +3.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.666"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "29.792.71"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008424"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.40"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "2.113.66"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.659"
$ws.Range("E23").Value = "  +2.97%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1516"
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.990"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.39"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.45"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.547"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.294"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.241"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05101"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7910"
$ws.Range("E34").Value = "  +3.86%  "
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.167"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").Value = "1.336.38"
$ws.Range("E38").Value = "  +9.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01888"
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9656"
$ws.Range("E41").Value = "  +6.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.058"
$ws.Range("E42").Value = "  +11.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.25"
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000127"
$ws.Range("E45").Value = "  +6.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.818"
$ws.Range("E46").Value = "  +3.82%  "
$ws.Range("D47").Value = "2.013.64"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.74"
$ws.Range("E48").Value = "  +3.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.799"
$ws.Range("E49").Value = "  +4.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5197"
$ws.Range("E50").Value = "  +0.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.046"
$ws.Range("E51").Value = "  +2.07%  "
